$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of 43 values for column A (rows 1-43)
$newValues = @(
    "old_DB_names",
    "Rv_ID",
    "gene_name",
    "2003A_Sassetti",
    "2003B_Sassetti",
    "2005_Rengarajan",
    "2006_Joshi_GI_1",
    "2006_Joshi_GI_2",
    "2011_Griffin_2",
    "2013_Zhang_1A",
    "2013_Zhang_1B",
    "2013_Zhang_2",
    "2013_Zhang_3A",
    "2013_Zhang_3B",
    "2013_Zhang_3C",
    "2013_Zhang_3D",
    "2016_Korte",
    "2018_Rittershaus_1A",
    "2018_Rittershaus_1B",
    "2012_Zhang",
    "2015_Kieser_GI_1",
    "2015_Kieser_GI_2",
    "2015_Kieser_GI_3",
    "2015_Mendum",
    "2016_Nambi",
    "2017_Mishra_1",
    "2017_Mishra_2",
    "2017_Xu_1A",
    "2017_Xu_1B",
    "2017_Xu_1C",
    "2017_Xu_1D",
    "2017_Xu_1E",
    "2018_Carey_1A",
    "2018_Carey_1B",
    "2018_Carey_1C",
    "2018_Carey_1D",
    "2018_Carey_1E",
    "2018_Carey_1F",
    "2018_Carey_1G",
    "2018_Carey_1H",
    "2017B_DeJesus_1A",
    "2017B_DeJesus_1B",
    "2017B_DeJesus_1C"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}

# Remove now-unused trailing rows (44-61) that held the old Rv gene name list
$ws.Range("A44:A61").EntireRow.Delete()

